$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Лист1" and name it "Лист2"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Лист2"

# --- populate Лист2 -------------------------------------------------
$ws2.Range("A1").Value = 5
$ws2.Range("B1").Value = "Коэффициент обратной связи"

$ws2.Range("B2").Value = "370 мВ"
$ws2.Range("B3").Value = "110 мВ"
$ws2.Range("C2").Value = "K_VT2"
$ws2.Range("C3").Value = "K_VT1"

$ws2.Range("A4").Value = 6
$ws2.Range("B4").Value = "Частотомер"

$ws2.Range("B5").Value = 1.35
$ws2.Range("C5").Value = "МГц"
$ws2.Range("D5").Value = 12
$ws2.Range("E5").Value = "В"

$ws2.Range("B6").Value = 1.346
$ws2.Range("C6").Value = "МГц"

$ws2.Range("B7").Value = 1.006
$ws2.Range("C7").Value = "МГц"
$ws2.Range("D7").Value = 11
$ws2.Range("E7").Value = "В"

$ws2.Range("B8").Value = 1.337
$ws2.Range("C8").Value = "МГц"
$ws2.Range("D8").Value = 12
$ws2.Range("E8").Value = "В"

$ws2.Range("B9").Value = 1.335
$ws2.Range("C9").Value = "МГц"
$ws2.Range("D9").Value = 11
$ws2.Range("E9").Value = "В"

$ws2.Range("B10").Value = 1.332
$ws2.Range("C10").Value = "МГц"
$ws2.Range("D10").Value = 3
$ws2.Range("E10").Value = "В"

$ws2.Range("B11").Value = "Коэффициент передачи при стандарных"
$ws2.Range("F11").Value = "Ku"
$ws2.Range("G11").Value = 52

$ws2.Range("A12").Value = 7
$ws2.Range("B12").Value = "Vп, В"
$ws2.Range("C12").Value = "Vамплитуды, мВ"

$ws2.Range("B13").Value = 12
$ws2.Range("C13").Value = 139

$ws2.Range("B14").Value = 11
$ws2.Range("C14").Value = 128

$ws2.Range("B15").Value = 10
$ws2.Range("C15").Value = 115

$ws2.Range("B16").Value = 9
$ws2.Range("C16").Value = 103

$ws2.Range("B17").Value = 8
$ws2.Range("C17").Value = 89

$ws2.Range("B18").Value = 7
$ws2.Range("C18").Value = 76

$ws2.Range("B19").Value = 6
$ws2.Range("C19").Value = 63

$ws2.Range("B20").Value = 5
$ws2.Range("C20").Value = 48

$ws2.Range("B21").Value = 4
$ws2.Range("C21").Value = 35

$ws2.Range("A22").Value = "R"
$ws2.Range("B22").Value = 734
$ws2.Range("C22").Value = "Ом"

$ws2.Range("A23").Value = "Q"
$ws2.Range("B23").Value = 24

$ws2.Range("A24").Value = "R2"
$ws2.Range("B24").Value = 530
$ws2.Range("C24").Value = "Ом"
$ws2.Range("D24").Value = "Q"
$ws2.Range("E24").Value = 150

$ws2.Range("A25").Value = "Уход частоты"

$ws2.Range("B26").Value = 1.0005
$ws2.Range("C26").Value = 12
$ws2.Range("D26").Value = "В"

$ws2.Range("B27").Value = 1.00068
$ws2.Range("C27").Value = 10
$ws2.Range("D27").Value = "В"

$ws2.Range("B28").Value = 989.8
$ws2.Range("C28").Value = 4
$ws2.Range("D28").Value = "В"

$ws2.Range("B29").Value = 989.6
$ws2.Range("C29").Value = 4

$ws2.Range("B30").Value = 971.7
$ws2.Range("C30").Value = 3

$ws2.Range("A31").Value = 17.9
$ws2.Range("B31").Value = "кГц"
$ws2.Range("C31").Value = 180
$ws2.Range("D31").Value = "Гц"

# column C on Лист2 gets a custom width (matches the source workbook's formatting)
$ws2.Columns.Item(3).ColumnWidth = 14.6

# make Лист2 the active sheet/tab and leave the cursor just below the data
$ws2.Activate() | Out-Null
$ws2.Range("A32").Select() | Out-Null
